$d = $word.ActiveDocument

# Locate the (currently empty) list paragraph that immediately follows the
# "git status" bullet - that's the paragraph the diff fills in with
# "git add sure.txt", followed by a brand-new bullet for the commit line.
$statusIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "git status") {
        $statusIndex = $i
        break
    }
}

if ($statusIndex -eq -1) {
    throw "Could not find the 'git status' paragraph"
}

$target = $d.Paragraphs.Item($statusIndex + 1)

# Fill the existing empty paragraph with the first new line.
$target.Range.InsertAfter("git add sure.txt")
$target = $d.Paragraphs.Item($statusIndex + 1)
$target.Range.Font.Size = 14
$target.Range.Font.SizeBi = 14

# Add a brand-new paragraph (matching list formatting) for the second line.
$target.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($statusIndex + 2)
$newPara.Range.InsertAfter("git commit -m `"This is our comments to commit)")
$newPara.Range.Font.Size = 14
$newPara.Range.Font.SizeBi = 14
